$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlContinuous = 1
$xlThin = 2
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10

function Set-Border($range, $edges) {
    foreach ($edge in $edges) {
        $range.Borders.Item($edge).LineStyle = $xlContinuous
        $range.Borders.Item($edge).Weight = $xlThin
    }
}

# --- Fill in the new roster data -----------------------------------------
# Column A first (top-to-bottom), then column B (top-to-bottom), so that the
# shared-string table ends up in the same order the source workbook has.
$ws.Range("A2").Value = "Sarah Lee"
$ws.Range("A3").Value = "David Smith"
$ws.Range("A4").Value = "Emily Clark"

$ws.Range("B2").Value = "Team 01"
$ws.Range("B3").Value = "Team 02"
$ws.Range("B4").Value = "Team 03"

# --- Fonts -----------------------------------------------------------------
# Header row grows from 16pt to 18pt, body rows grow from 12pt to 14pt.
$ws.Range("A1:B1").Font.Size = 18
$ws.Range("A2:D5").Font.Size = 14

# --- Borders -----------------------------------------------------------------
# Column A gets a left/top/bottom "box" (no right edge - it abuts column B).
Set-Border $ws.Range("A1") @($xlEdgeLeft, $xlEdgeTop, $xlEdgeBottom)
Set-Border $ws.Range("A2") @($xlEdgeLeft, $xlEdgeTop)
Set-Border $ws.Range("A3") @($xlEdgeLeft)
Set-Border $ws.Range("A4") @($xlEdgeLeft, $xlEdgeBottom)

# Column B gets a full box around each cell (it's the right-most visible column).
Set-Border $ws.Range("B1") @($xlEdgeLeft, $xlEdgeTop, $xlEdgeBottom, $xlEdgeRight)
Set-Border $ws.Range("B2") @($xlEdgeLeft, $xlEdgeTop, $xlEdgeRight)
Set-Border $ws.Range("B3") @($xlEdgeLeft, $xlEdgeRight)
Set-Border $ws.Range("B4") @($xlEdgeLeft, $xlEdgeBottom, $xlEdgeRight)

# --- Extend the used range down to row 5 / out to column D -----------------
$ws.Range("C2:D5").Font.Size = 14
$ws.Range("A5:D5").Font.Size = 14

# --- Column widths / row heights -------------------------------------------
$ws.Columns("A").ColumnWidth = 19.3
$ws.Columns("B").ColumnWidth = 24.6

$ws.Rows("1").RowHeight = 23.25
$ws.Rows("2").RowHeight = 18.75
$ws.Rows("3").RowHeight = 18.75
$ws.Rows("4").RowHeight = 18.75
$ws.Rows("5").RowHeight = 18.75

# --- Selection ---------------------------------------------------------------
$ws.Range("D11").Select() | Out-Null
